$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "t27TflGL"
$ws.Range("B2").Value = "trashboatsr"
$ws.Range("C2").Value = 1890
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = "https://lichess.org/t27TflGL"
$ws.Range("F2").Value = 4159
$ws.Range("G2").Value = "NO"
$ws.Range("H2").Value = "blank"
$ws.Range("I2").Value = "NO"
